# Auto-generated edit script applying the cryptos.xlsx diff
# (refreshed coin prices / 1h volume %, and two row-pair swaps: rows 26<->27, 44<->45)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.506.41'
$ws.Range('E2').Value = '  +1.76%  '
$ws.Range('D3').Value = '2.656.66'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '581.71'
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.83'
$ws.Range('E6').Value = '  +1.09%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.603'
$ws.Range('E8').Value = '  +1.09%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.61'
$ws.Range('E9').Value = '  +0.85%  '
$ws.Range('E10').Value = '  +4.07%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.383'
$ws.Range('E11').Value = '  +2.92%  '
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('D13').Value = '3.123.61'
$ws.Range('E13').Value = '  +2.14%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.12'
$ws.Range('E14').Value = '  +6.46%  '
$ws.Range('D15').Value = '61.396.14'
$ws.Range('E15').Value = '  +1.61%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000147'
$ws.Range('E16').Value = '  +3.84%  '
$ws.Range('D17').Value = '2.662.11'
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.69'
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('E19').Value = '  +2.98%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '355.66'
$ws.Range('E20').Value = '  +2.18%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.92'
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.526'
$ws.Range('E23').Value = '  +0.62%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '64.75'
$ws.Range('E24').Value = '  +2.41%  '
$ws.Range('E25').Value = '  +2.87%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.50'
$ws.Range('E26').Value = '  +5.13%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.996'
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.01'
$ws.Range('E28').Value = '  +7.99%  '
$ws.Range('D29').Value = '0.0₃0824'
$ws.Range('E29').Value = '  +3.22%  '
$ws.Range('E30').Value = '  +8.31%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '169.61'
$ws.Range('E31').Value = '  +3.17%  '
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.19'
$ws.Range('E33').Value = '  +3.70%  '
$ws.Range('E34').Value = '  +15.03%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.69'
$ws.Range('E35').Value = '  +8.56%  '
$ws.Range('E36').Value = '  +8.56%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +19.19%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.73'
$ws.Range('E38').Value = '  +5.16%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '344.59'
$ws.Range('E39').Value = '  +9.66%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.15'
$ws.Range('E40').Value = '  +6.36%  '
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.42'
$ws.Range('E42').Value = '  +7.08%  '
$ws.Range('E43').Value = '  +5.17%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '21.22'
$ws.Range('E44').Value = '  +5.34%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '20.70'
$ws.Range('E45').Value = '  +4.20%  '
$ws.Range('E46').Value = '  +5.16%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '135.79'
$ws.Range('E47').Value = '  +0.54%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.629'
$ws.Range('E48').Value = '  +3.61%  '
$ws.Range('E49').Value = '  +1.11%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.997'
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('D51').Value = '2.103.11'
$ws.Range('E51').Value = '  +3.50%  '
